# Re-run SGNN to annotate dialog acts following clean up work to the original transcripts.
# Updates the DAMSLTag (column I) and DialogAct (column J) values for the rows
# whose automated annotation changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 4;   Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 11;  Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 14;  Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 17;  Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 19;  Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 26;  Tag = "ba"; Act = "Appreciation" },
    @{ Row = 41;  Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 42;  Tag = "qy"; Act = "Yes-No-Question" },
    @{ Row = 43;  Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 55;  Tag = "qy"; Act = "Yes-No-Question" },
    @{ Row = 57;  Tag = "qy"; Act = "Yes-No-Question" },
    @{ Row = 59;  Tag = "b";  Act = "Acknowledge (Backchannel)" },
    @{ Row = 64;  Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 78;  Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 93;  Tag = "ba"; Act = "Appreciation" },
    @{ Row = 98;  Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 106; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 107; Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 109; Tag = "ba"; Act = "Appreciation" },
    @{ Row = 115; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 120; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 169; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 180; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 187; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 195; Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 204; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 206; Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 211; Tag = "sd"; Act = "Statement-non-opinion" }
)

foreach ($u in $updates) {
    $ws.Range("I" + $u.Row).Value = $u.Tag
    $ws.Range("J" + $u.Row).Value = $u.Act
}
